# Update "想去人数" (number of people interested) counts that changed
# between the two scrape runs, on both the "展览" sheet and the
# "全部类型" aggregate sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 1065
$wsExpo.Range("F6").Value = 14019
$wsExpo.Range("F7").Value = 15335
$wsExpo.Range("F23").Value = 5961
$wsExpo.Range("F29").Value = 94
$wsExpo.Range("F30").Value = 435

# --- Sheet "全部类型" (all categories combined) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 1065
$wsAll.Range("F7").Value = 14019
$wsAll.Range("F8").Value = 15335
$wsAll.Range("F25").Value = 5961
$wsAll.Range("F31").Value = 94
$wsAll.Range("F32").Value = 435
